# Update "想去人数" (want-to-go count) figures in column F across the
# "展览" (rId1), "演出" (rId2) and "全部类型" (rId4) worksheets.
# "本地生活" (rId3) is left untouched.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        3  = 677
        4  = 114
        5  = 2110
        7  = 10976
        8  = 187
        9  = 165
        10 = 298
        11 = 211
        12 = 10817
        13 = 437
        14 = 1129
        15 = 23
        16 = 753
        17 = 5419
        18 = 80
        19 = 3401
    }
    "演出" = @{
        3 = 564
    }
    "全部类型" = @{
        3  = 677
        5  = 114
        6  = 2110
        7  = 564
        10 = 10976
        11 = 187
        12 = 165
        13 = 298
        14 = 211
        15 = 10818
        16 = 437
        17 = 1129
        18 = 23
        19 = 753
        20 = 5419
        21 = 80
        22 = 3401
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
